# Append new booking rows (37-44) to the bookings sheet, as recorded by the
# Snow Liwa booking app on 2025-12-04.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A booking_id, B created_at, C name, D phone, E tickets,
#          F ticket_price, G total_amount, H status, I payment_intent_id,
#          J payment_status, K redirect_url, L notes

# Phone-like / numeric-looking text columns must be forced to text so that
# leading zeros (and the purely-numeric "1234" name) survive the round trip.
$ws.Range("D37:D38").NumberFormat = "@"
$ws.Range("D40:D44").NumberFormat = "@"
$ws.Range("C44").NumberFormat = "@"

# Row 37
$ws.Range("A37").Value = "SL-20251204-001"
$ws.Range("B37").Value = "2025-12-04 00:48:53"
$ws.Range("C37").Value = "fahaf"
$ws.Range("D37").Value = "0502992692"
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 175
$ws.Range("G37").Value = 175
$ws.Range("H37").Value = "pending"
$ws.Range("J37").Value = "pending"

# Row 38
$ws.Range("A38").Value = "SL-20251204-002"
$ws.Range("B38").Value = "2025-12-04 00:49:18"
$ws.Range("C38").Value = "fahaf"
$ws.Range("D38").Value = "0502992692"
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = 175
$ws.Range("G38").Value = 175
$ws.Range("H38").Value = "pending"
$ws.Range("J38").Value = "pending"

# Row 39
$ws.Range("A39").Value = "SL-20251204-003"
$ws.Range("B39").Value = "2025-12-04 02:21:57"
$ws.Range("C39").Value = "كل أيام الأسبوع"
$ws.Range("D39").Value = "4:00pm - 12:00am"
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = 175
$ws.Range("G39").Value = 175
$ws.Range("H39").Value = "pending"
$ws.Range("J39").Value = "pending"

# Row 40
$ws.Range("A40").Value = "SL-20251204-004"
$ws.Range("B40").Value = "2025-12-04 04:27:13"
$ws.Range("C40").Value = "fahad"
$ws.Range("D40").Value = "0502992932"
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 175
$ws.Range("G40").Value = 175
$ws.Range("H40").Value = "pending"
$ws.Range("J40").Value = "pending"

# Row 41
$ws.Range("A41").Value = "SL-20251204-005"
$ws.Range("B41").Value = "2025-12-04 04:27:13"
$ws.Range("C41").Value = "fahad"
$ws.Range("D41").Value = "0502992932"
$ws.Range("E41").Value = 1
$ws.Range("F41").Value = 175
$ws.Range("G41").Value = 175
$ws.Range("H41").Value = "pending"
$ws.Range("J41").Value = "pending"

# Row 42
$ws.Range("A42").Value = "SL-20251204-006"
$ws.Range("B42").Value = "2025-12-04 05:03:20"
$ws.Range("C42").Value = "Fahad"
$ws.Range("D42").Value = "0502992932"
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 175
$ws.Range("G42").Value = 175
$ws.Range("H42").Value = "pending"
$ws.Range("J42").Value = "pending"

# Row 43
$ws.Range("A43").Value = "SL-20251204-007"
$ws.Range("B43").Value = "2025-12-04 05:04:26"
$ws.Range("C43").Value = "Fahad"
$ws.Range("D43").Value = "0502992932"
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = 175
$ws.Range("G43").Value = 175
$ws.Range("H43").Value = "pending"
$ws.Range("J43").Value = "pending"

# Row 44
$ws.Range("A44").Value = "SL-20251204-008"
$ws.Range("B44").Value = "2025-12-04 06:15:53"
$ws.Range("C44").Value = "1234"
$ws.Range("D44").Value = "0502992932"
$ws.Range("E44").Value = 1
$ws.Range("F44").Value = 175
$ws.Range("G44").Value = 175
$ws.Range("H44").Value = "pending"
$ws.Range("J44").Value = "pending"
